$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the last-changed date serial for each row.
# Update all data rows (2 through 45) from 45204 to 45205 (2023-10-05 -> 2023-10-06).
for ($row = 2; $row -le 45; $row++) {
    $ws.Cells.Item($row, 3).Value2 = 45205
}
